$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - new observation record (Spillkråka / Dryocopus martius)
$ws.Cells.Item(8, 1).Value = 131073777      # A8 Id
$ws.Cells.Item(8, 2).Value = 57881          # B8 Taxonsorteringsordning
$ws.Cells.Item(8, 4).Value = "NT"           # D8 Rödlistade
$ws.Cells.Item(8, 5).Value = 100049         # E8 TaxonId
$ws.Cells.Item(8, 6).Value = "Spillkråka"   # F8 Artnamn
$ws.Cells.Item(8, 7).Value = "Dryocopus martius"   # G8 Vetenskapligt namn
$ws.Cells.Item(8, 8).Value = "(Linnaeus, 1758)"    # H8 Auktor
$ws.Cells.Item(8, 13).Value = "färska spår"        # M8 Aktivitet
$ws.Cells.Item(8, 16).Value = "Stötetorpet, Stötetorpet, Boh"   # P8 Lokalnamn
$ws.Cells.Item(8, 17).Value = 311113        # Q8 Ost
$ws.Cells.Item(8, 18).Value = 6410575       # R8 Nord
$ws.Cells.Item(8, 19).Value = 10            # S8 Noggrannhet
$ws.Cells.Item(8, 20).Value = "Västra Götaland"    # T8 Län
$ws.Cells.Item(8, 21).Value = "Kungälv"     # U8 Kommun
$ws.Cells.Item(8, 22).Value = "Bohuslän"    # V8 Provins
$ws.Cells.Item(8, 23).Value = "Harestad"    # W8 Socken

# Startdatum/Slutdatum are stored as plain text ("YYYY-MM-DD"), not real dates,
# throughout this workbook. Force text formatting on these two cells so the
# values round-trip as strings instead of being auto-parsed into date serials.
$ws.Cells.Item(8, 25).NumberFormat = "@"
$ws.Cells.Item(8, 25).Value = "2026-02-08"  # Y8 Startdatum
$ws.Cells.Item(8, 27).NumberFormat = "@"
$ws.Cells.Item(8, 27).Value = "2026-02-08"  # AA8 Slutdatum

$ws.Cells.Item(8, 30).Value = $false        # AD8 Ej återfunnen
$ws.Cells.Item(8, 31).Value = $false        # AE8 Osäker artbestämning
$ws.Cells.Item(8, 33).Value = $false        # AG8 Ospontan

$ws.Cells.Item(8, 49).Value = "Linus Lundin"   # AW8 Rapportör
$ws.Cells.Item(8, 50).Value = "Linus Lundin"   # AX8 Observatörer
